# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" positioned before the existing "2022-Q3"
#    sheet (i.e. right after "总计").
# 2) Populate the new sheet with the fund holding table for 2022-Q4, matching
#    the look/feel (header row + index column styling) of the other quarter
#    sheets.
# 3) Update the "总计" (summary) sheet so row 2 becomes the 2022-Q4 entry and
#    every other quarter's row shifts down by one (2021-Q4 becomes row 6).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create + place the new sheet
# ---------------------------------------------------------------------------
# NOTE: after a structural change like Worksheets.Add(), any previously held
# worksheet variable can end up re-resolving to the NEW sheet instead of the
# original one, so every reference used below is re-fetched fresh by name
# right when it's needed.
$insertBeforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($insertBeforeSheet)
$q4Sheet.Name = "2022-Q4"

# Re-fetch both sheets by name now that the workbook structure changed.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Item("2022-Q4")

# ---------------------------------------------------------------------------
# Step 2: formatting — copy header-row style and index-column style from the
# neighbouring "2022-Q3" sheet so the new sheet matches the existing look.
# ---------------------------------------------------------------------------
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

$q3Sheet.Range("A2").Copy()
$q4Sheet.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2b: header labels
# ---------------------------------------------------------------------------
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# Step 2c: data rows (A = numeric index, B..G = text, H = numeric rank)
# ---------------------------------------------------------------------------
$q4Data = @(
    @{ A=0; B="002601"; C="中银证券价值精选灵活配置混合"; D="5.15"; E="92.90"; F="3.89"; G="0.2003"; H=8 },
    @{ A=1; B="001543"; C="宝盈新锐灵活配置混合A";       D="2.28"; E="91.59"; F="5.22"; G="0.1190"; H=4 },
    @{ A=2; B="006323"; C="合煦智远嘉选混合A";           D="0.70"; E="73.98"; F="5.59"; G="0.0391"; H=2 },
    @{ A=3; B="007578"; C="宝盈新锐灵活配置混合C";       D="0.21"; E="91.59"; F="5.22"; G="0.0110"; H=4 },
    @{ A=4; B="006324"; C="合煦智远嘉选混合C";           D="0.14"; E="73.98"; F="5.59"; G="0.0078"; H=2 }
)

# Text-typed columns (B..G) get a leading apostrophe so the numeric-looking
# strings ("5.15", "002601", ...) are stored as text instead of being
# auto-coerced to numbers — then the auto-applied "quote prefix" cell style
# is reset back to Normal so no stray NumberFormat/style sticks around
# (matches the plain, un-styled data cells in the other quarter sheets).
$row = 2
foreach ($rec in $q4Data) {
    $q4Sheet.Range("A$row").Value = $rec.A

    $q4Sheet.Range("B$row").Value = "'" + $rec.B
    $q4Sheet.Range("C$row").Value = "'" + $rec.C
    $q4Sheet.Range("D$row").Value = "'" + $rec.D
    $q4Sheet.Range("E$row").Value = "'" + $rec.E
    $q4Sheet.Range("F$row").Value = "'" + $rec.F
    $q4Sheet.Range("G$row").Value = "'" + $rec.G
    $q4Sheet.Range("B$row`:G$row").Style = "Normal"

    $q4Sheet.Range("H$row").Value = $rec.H
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Step 3: update the "总计" summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryData = @(
    @{ Row=2; A=0; B="2022-Q4"; C=5; D=0.38 },
    @{ Row=3; A=1; B="2022-Q3"; C=6; D=0.17 },
    @{ Row=4; A=2; B="2022-Q2"; C=2; D=0.03 },
    @{ Row=5; A=3; B="2022-Q1"; C=4; D=0.03 },
    @{ Row=6; A=4; B="2021-Q4"; C=5; D=1.16 }
)

foreach ($rec in $summaryData) {
    $r = $rec.Row
    $summary.Range("A$r").Value = $rec.A
    $summary.Range("B$r").Value = "'" + $rec.B
    $summary.Range("B$r").Style = "Normal"
    $summary.Range("C$r").Value = $rec.C
    $summary.Range("D$r").Value = $rec.D
}

Write-Host "2022-Q4 sheet + 总计 summary updated"
